# Updated cryptos list on Tue Oct 24 16:54:15 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns for each coin row.
# D-column values are decimal-looking text (e.g. "225.00", "33.998.40") that
# must stay literal strings -- Excel would otherwise auto-coerce a plain
# `.Value = "225.00"` assignment into the number 225, dropping the trailing
# zero / thousands-style dots. Forcing the cell to Text format first
# preserves the exact characters, and resetting the style back to "Normal"
# afterwards avoids leaving a stray number-format behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceCell($row, $value) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

function Set-VolumeCell($row, $value) {
    $ws.Range("E$row").Value = $value
}

# row -> (new price text or $null if unchanged), new volume text
$updates = @(
    @(2,  "33.998.40",  "  +9.46%  "),
    @(3,  "1.785.81",   "  +6.20%  "),
    @(4,  "0.999",      "  +0.21%  "),
    @(5,  "225.00",     "  +2.16%  "),
    @(6,  "0.557",      "  +4.34%  "),
    @(7,  "0.999",      "  +0.13%  "),
    @(8,  "30.82",      "  +5.13%  "),
    @(9,  "46.32",      "  +3.42%  "),
    @(10, "0.278",      "  +4.70%  "),
    @(11, $null,        "  +3.41%  "),
    @(12, "0.0925",     "  +1.91%  "),
    @(13, "2.044.59",   "  +6.40%  "),
    @(14, "1.794.80",   "  +6.68%  "),
    @(15, "0.627",      "  +3.42%  "),
    @(16, "33.942.06",  "  +9.51%  "),
    @(17, $null,        "  -2.52%  "),
    @(18, $null,        "  +2.53%  "),
    @(19, "68.51",      "  +3.37%  "),
    @(20, "251.76",     "  +1.83%  "),
    @(21, $null,        "  +2.74%  "),
    @(22, "0.998",      "  +0.03%  "),
    @(23, "10.28",      "  +2.90%  "),
    @(24, $null,        "  -0.96%  "),
    @(25, $null,        "  -0.06%  "),
    @(26, "158.27",     "  -0.41%  "),
    @(27, "16.46",      "  +3.59%  "),
    @(28, $null,        "  +2.02%  "),
    @(29, $null,        "  +3.64%  "),
    @(30, "0.999",      "  +0.17%  "),
    @(31, $null,        "  +8.20%  "),
    @(32, "0.0512",     "  +3.26%  "),
    @(33, $null,        "  +4.06%  "),
    @(34, $null,        "  +5.82%  "),
    @(35, "1.498.11",   "  -1.15%  "),
    @(36, $null,        "  +1.30%  "),
    @(37, $null,        "  +3.54%  "),
    @(38, "0.628",      "  +2.94%  "),
    @(39, "0.0185",     "  +3.03%  "),
    @(40, "83.16",      "  -1.52%  "),
    @(41, $null,        "  +2.71%  "),
    @(42, $null,        "  +1.32%  "),
    @(43, "0.886",      "  +5.48%  "),
    @(44, $null,        "  +1.08%  "),
    @(45, "0.0511",     "  +1.51%  "),
    @(46, $null,        "  +3.77%  "),
    @(47, "1.939.46",   "  +6.88%  "),
    @(48, "5.74",       "  +2.76%  "),
    @(49, "1.00",       "  +0.03%  "),
    @(50, "11.93",      "  +12.13%  "),
    @(51, "50.82",      "  -1.27%  ")
)

foreach ($u in $updates) {
    $row = $u[0]
    $price = $u[1]
    $volume = $u[2]

    if ($null -ne $price) {
        Set-PriceCell $row $price
    }
    Set-VolumeCell $row $volume
}
